# [FIX] filter charge_type (expense, gl, revenue and trial balance web)
#
# The "Project Report" filter block (rows ~2-19) lists one filter label per
# row in column A. This adds a "Charge Type" filter line right above
# "Run By" / "Run Date" (previously "Charge Type" only showed up as a
# column header further down the sheet). Inserting the row pushes the
# blank spacer row and the big header row (with one cell per report
# column) down by one, which is also reflected in the new sheet
# dimension (BJ21 -> BJ22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current "Run By" row (row 20) - this shifts
# the blank spacer row (old 20) to 21 and the column-header row (old 21)
# to 22, carrying their existing content/formatting along automatically
# (the new row inherits the look of the row above it, i.e. the same
# label style used by the other filter-label cells in column A).
$ws.Rows.Item(20).Insert()

# Populate the new filter row with the "Charge Type" label.
$ws.Range("A20").Value = "Charge Type"
